$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.404.90"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "1.872.27"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.32"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4713"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2869"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06494"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.87"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "100.33"
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07801"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.871.65"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7281"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.170"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.03"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "30.384.58"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.10"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007484"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").Value = "2.115.71"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.328"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.333"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.043"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.33"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.96"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.895"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09675"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.321"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.490"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.230"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.152"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04805"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.126"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6897"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.836"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.57"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.295"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.957"
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4221"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8248"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.90"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.705"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.016"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.97"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "883.69"
$ws.Range("E51").Value = "  -3.87%  "